$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 891 (shifts existing data down by 4, producing the
# dimension growth from R962 to R966 and re-numbering every later row automatically).
$ws.Rows("891:894").Insert()

# New row 891: week of 44461, variety Conconina(o)
$ws.Range("A891").Value = 3
$ws.Range("B891").Value = 'Femacal de La Calera'
$ws.Range("C891").Value = 'Coquimbo'
$ws.Range("D891").Value = 44461
$ws.Range("E891").Value = 5
$ws.Range("F891").Value = 100112033
$ws.Range("G891").Value = 'Lechuga'
$ws.Range("H891").Value = 'Conconina(o)'
$ws.Range("I891").Value = 'Primera'
$ws.Range("J891").Value = 130
$ws.Range("K891").Value = 4500
$ws.Range("L891").Value = 5000
$ws.Range("M891").Value = 4731
$ws.Range("N891").Value = '$/caja 10 unidades'
$ws.Range("O891").Value = 'Provincia de Quillota'
$ws.Range("P891").Value = 473
$ws.Range("Q891").Value = 10
$ws.Range("R891").Value = 'Hortaliza'

# New row 892: week of 44461, variety Escarola
$ws.Range("A892").Value = 3
$ws.Range("B892").Value = 'Femacal de La Calera'
$ws.Range("C892").Value = 'Coquimbo'
$ws.Range("D892").Value = 44461
$ws.Range("E892").Value = 5
$ws.Range("F892").Value = 100112033
$ws.Range("G892").Value = 'Lechuga'
$ws.Range("H892").Value = 'Escarola'
$ws.Range("I892").Value = 'Primera'
$ws.Range("J892").Value = 115
$ws.Range("K892").Value = 5500
$ws.Range("L892").Value = 6000
$ws.Range("M892").Value = 5739
$ws.Range("N892").Value = '$/caja 15 unidades'
$ws.Range("O892").Value = 'Provincia de Quillota'
$ws.Range("P892").Value = 383
$ws.Range("Q892").Value = 15
$ws.Range("R892").Value = 'Hortaliza'

# New row 893: week of 44461, variety Francesa morada
$ws.Range("A893").Value = 3
$ws.Range("B893").Value = 'Femacal de La Calera'
$ws.Range("C893").Value = 'Coquimbo'
$ws.Range("D893").Value = 44461
$ws.Range("E893").Value = 5
$ws.Range("F893").Value = 100112033
$ws.Range("G893").Value = 'Lechuga'
$ws.Range("H893").Value = 'Francesa morada'
$ws.Range("I893").Value = 'Primera'
$ws.Range("J893").Value = 120
$ws.Range("K893").Value = 4500
$ws.Range("L893").Value = 4800
$ws.Range("M893").Value = 4650
$ws.Range("N893").Value = '$/caja 18 unidades'
$ws.Range("O893").Value = 'Provincia de Quillota'
$ws.Range("P893").Value = 258
$ws.Range("Q893").Value = 18
$ws.Range("R893").Value = 'Hortaliza'

# New row 894: week of 44461, variety Marina
$ws.Range("A894").Value = 3
$ws.Range("B894").Value = 'Femacal de La Calera'
$ws.Range("C894").Value = 'Coquimbo'
$ws.Range("D894").Value = 44461
$ws.Range("E894").Value = 5
$ws.Range("F894").Value = 100112033
$ws.Range("G894").Value = 'Lechuga'
$ws.Range("H894").Value = 'Marina'
$ws.Range("I894").Value = 'Primera'
$ws.Range("J894").Value = 60
$ws.Range("K894").Value = 4500
$ws.Range("L894").Value = 4500
$ws.Range("M894").Value = 4500
$ws.Range("N894").Value = '$/caja 18 unidades'
$ws.Range("O894").Value = 'Provincia de Quillota'
$ws.Range("P894").Value = 250
$ws.Range("Q894").Value = 18
$ws.Range("R894").Value = 'Hortaliza'
